$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text for the quiz instructions in A3
$ws.Range("A3").Value = "You'll be given a quiz before the story for the 4 questions associated with the assigned perspective. "

# Update text for the "must pass this quiz" instructions in A4
$ws.Range("A4").Value = "You must pass this quiz to move onto the story. "

# Adjust row heights for rows 3 and 4
$ws.Rows(3).RowHeight = 109.55
$ws.Rows(4).RowHeight = 64.55
